# Applies the "Historico Diseno.docx" edit:
#  1) Collapses the "Se cambia ... pagina ... usuarios" sentence (previously
#     split across runs around two <w:proofErr> gramStart/gramEnd markers)
#     into a single run.
#  2) Collapses the "Se pasan los scripts ... archivos externos" sentence
#     (previously split around a <w:proofErr> gramStart/gramEnd pair) into a
#     single run.
#  3) Appends two new log paragraphs ("Saque titulo pagina ..." and "Sigue
#     el problema del footer ...") plus one extra blank paragraph right
#     after the "V6.1  Agregue icono en index" entry.

$d = $word.ActiveDocument

# --- Change 1: merge the "Se cambia el pagina ... usuarios" run split ---
$text1 = "Se cambia el página “quienes somos” por la página “usuarios”"
$d.Content.Find.Execute($text1, $true, $false, $false, $false, $false, $true, 1, $false, $text1, 2) | Out-Null

# --- Change 2: merge the "Se pasan los scripts ... carpeta" run split ---
$text2 = "Se pasan los scripts de los formularios  a archivos externos en la carpeta "
$d.Content.Find.Execute($text2, $true, $false, $false, $false, $false, $true, 1, $false, $text2, 2) | Out-Null

# --- Change 3: insert the two new paragraphs + one blank paragraph ---
# Locate the "V6.1  Agregue icono en index" paragraph as an anchor (rather
# than a hard-coded paragraph index) so the edit is resilient to any minor
# paragraph-numbering differences.
$anchorRng = $d.Content
$anchorRng.Find.Execute("Agregue icono en", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$anchorPara = $anchorRng.Paragraphs(1)

# Create a placeholder paragraph right after the anchor, then replace its
# contents with the two fully-formed new paragraphs (keeps the required
# <w:proofErr> spell-check markers intact).
$anchorPara.Range.InsertParagraphAfter() | Out-Null
$newParas = $anchorPara.Next()
$newParas.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:tab/><w:t xml:space="preserve">Saque titulo </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>pagina</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> que se superponía con el menú en Usuarios, contactos, recomendaciones, juegos y consolas. En su lugar puse un párrafo centrado en cada una de ella en el </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>body</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>, hay que cambiar el color del menú, cuesta leerlo contra el fondo de la pantalla y agrandar el titulo con el mismo color del menú porque tampoco se lee.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Sigue el problema del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>footer</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> en las </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>paginas</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> de celulares</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null

# Locate the just-inserted "Sigue el problema del footer ..." paragraph and
# append one more, genuinely empty, paragraph after it.
$footerRng = $d.Content
$footerRng.Find.Execute("Sigue el problema del", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$footerPara = $footerRng.Paragraphs(1)
$footerPara.Range.InsertParagraphAfter() | Out-Null
$blankPara = $footerPara.Next()
$blankPara.Range.InsertXML('<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>') | Out-Null
